$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2195121951219512
$ws.Range("C2").Value = 0.5040650406504065
$ws.Range("J2").Value = 0.01084010840108401
$ws.Range("P2").Value = 0.1653116531165312
$ws.Range("S2").Value = 0.1002710027100271
$ws.Range("B3").Value = 0.005291005291005291
$ws.Range("C3").Value = 0.02116402116402116
$ws.Range("J3").Value = 0.02645502645502645
$ws.Range("P3").Value = 0.7724867724867724
$ws.Range("S3").Value = 0.1746031746031746
$ws.Range("P4").Value = 0.6268656716417911
$ws.Range("S4").Value = 0.373134328358209
$ws.Range("B6").Value = 0.05752212389380531
$ws.Range("D6").Value = 0.01327433628318584
$ws.Range("F6").Value = 0.04424778761061947
$ws.Range("J6").Value = 0.2743362831858407
$ws.Range("O6").Value = 0.03097345132743363
$ws.Range("Q6").Value = 0.163716814159292
$ws.Range("R6").Value = 0.05309734513274336
$ws.Range("S6").Value = 0.3628318584070797
$ws.Range("B7").Value = 0.1004016064257028
$ws.Range("D7").Value = 0.02008032128514056
$ws.Range("F7").Value = 0.06827309236947791
$ws.Range("J7").Value = 0.1164658634538153
$ws.Range("O7").Value = 0.04016064257028112
$ws.Range("Q7").Value = 0.1927710843373494
$ws.Range("R7").Value = 0.04016064257028112
$ws.Range("S7").Value = 0.4216867469879518
$ws.Range("B8").Value = 0.1037924151696607
$ws.Range("D8").Value = 0.02794411177644711
$ws.Range("F8").Value = 0.05389221556886228
$ws.Range("J8").Value = 0.1037924151696607
$ws.Range("O8").Value = 0.01397205588822355
$ws.Range("Q8").Value = 0.1856287425149701
$ws.Range("R8").Value = 0.08582834331337326
$ws.Range("S8").Value = 0.4251497005988024
$ws.Range("B9").Value = 0.1214285714285714
$ws.Range("D9").Value = 0.05714285714285714
$ws.Range("F9").Value = 0.1
$ws.Range("J9").Value = 0.07857142857142857
$ws.Range("O9").Value = 0.007142857142857143
$ws.Range("Q9").Value = 0.2214285714285714
$ws.Range("R9").Value = 0.05714285714285714
$ws.Range("S9").Value = 0.3571428571428572
$ws.Range("B10").Value = 0.1306607275426875
$ws.Range("D10").Value = 0.02821083890126206
$ws.Range("F10").Value = 0.066815144766147
$ws.Range("J10").Value = 0.1121009651076466
$ws.Range("O10").Value = 0.02301410541945063
$ws.Range("Q10").Value = 0.2093541202672606
$ws.Range("R10").Value = 0.06607275426874536
$ws.Range("S10").Value = 0.3637713437268003
$ws.Range("G11").Value = 0.1432584269662921
$ws.Range("J11").Value = 0.08146067415730338
$ws.Range("K11").Value = 0.1797752808988764
$ws.Range("L11").Value = 0.5814606741573034
$ws.Range("S11").Value = 0.01404494382022472
$ws.Range("G12").Value = 0.8056872037914692
$ws.Range("J12").Value = 0.1706161137440758
$ws.Range("K12").Value = 0.004739336492890996
$ws.Range("L12").Value = 0.004739336492890996
$ws.Range("S12").Value = 0.01421800947867299
$ws.Range("G13").Value = 0.5573770491803278
$ws.Range("J13").Value = 0.3278688524590164
$ws.Range("S13").Value = 0.1147540983606557
$ws.Range("F15").Value = 0.02690582959641256
$ws.Range("H15").Value = 0.2242152466367713
$ws.Range("I15").Value = 0.08520179372197309
$ws.Range("J15").Value = 0.304932735426009
$ws.Range("K15").Value = 0.07174887892376682
$ws.Range("M15").Value = 0.01345291479820628
$ws.Range("O15").Value = 0.02690582959641256
$ws.Range("S15").Value = 0.2466367713004484
$ws.Range("F16").Value = 0.02074688796680498
$ws.Range("H16").Value = 0.1784232365145228
$ws.Range("I16").Value = 0.04149377593360996
$ws.Range("J16").Value = 0.4522821576763486
$ws.Range("K16").Value = 0.1161825726141079
$ws.Range("M16").Value = 0.01659751037344398
$ws.Range("N16").Value = 0.004149377593360996
$ws.Range("O16").Value = 0.05394190871369295
$ws.Range("S16").Value = 0.1161825726141079
$ws.Range("F17").Value = 0.02024291497975709
$ws.Range("H17").Value = 0.2044534412955465
$ws.Range("I17").Value = 0.05870445344129555
$ws.Range("J17").Value = 0.4170040485829959
$ws.Range("K17").Value = 0.1214574898785425
$ws.Range("M17").Value = 0.02024291497975709
$ws.Range("O17").Value = 0.04453441295546558
$ws.Range("S17").Value = 0.1133603238866397
$ws.Range("F18").Value = 0.006329113924050633
$ws.Range("H18").Value = 0.1645569620253164
$ws.Range("I18").Value = 0.06329113924050633
$ws.Range("J18").Value = 0.4620253164556962
$ws.Range("K18").Value = 0.08860759493670886
$ws.Range("M18").Value = 0.03164556962025317
$ws.Range("O18").Value = 0.05063291139240506
$ws.Range("S18").Value = 0.1329113924050633
$ws.Range("F19").Value = 0.01795063575168287
$ws.Range("H19").Value = 0.212415856394914
$ws.Range("I19").Value = 0.05534779356768885
$ws.Range("J19").Value = 0.3754674644727001
$ws.Range("K19").Value = 0.1271503365744203
$ws.Range("M19").Value = 0.02991772625280479
$ws.Range("N19").Value = 0.0007479431563201197
$ws.Range("O19").Value = 0.06955871353777113
$ws.Range("S19").Value = 0.1114435302916978
